$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 135.71428
$ws.Range("J2").Value = 563.3333
$ws.Range("H2").Value = 264
$ws.Range("M2").Value = -22.71428
$ws.Range("L2").Value = 563.3333
$ws.Range("K2").Value = 135.71428
$ws.Range("N2").Value = -789.3333
$ws.Range("K5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("I19").Value = 270.85715
$ws.Range("H19").Value = 338.19232
$ws.Range("M19").Value = -95.85714999999999
$ws.Range("K19").Value = 270.85715
$ws.Range("M32").Value = -1374
$ws.Range("N32").Value = -1612.4
$ws.Range("L32").Value = 960.4
$ws.Range("K32").Value = 1700
$ws.Range("J32").Value = 960.4
$ws.Range("I32").Value = 1700
$ws.Range("H32").Value = 1083.6666
$ws.Range("K64").Value = 166668960
$ws.Range("N64").Value = -5072
$ws.Range("J64").Value = 4576
$ws.Range("I64").Value = 166668960
$ws.Range("H64").Value = 62503724
$ws.Range("M64").Value = -166668712
$ws.Range("L64").Value = 4576
$ws.Range("N67").Value = -6292
$ws.Range("I67").Value = 166668960
$ws.Range("J67").Value = 4576
$ws.Range("H67").Value = 62503724
$ws.Range("M67").Value = -166668102
$ws.Range("L67").Value = 4576
$ws.Range("K67").Value = 166668960
$ws.Range("K74").Value = 2513.8
$ws.Range("J74").Value = 2866.6667
$ws.Range("I74").Value = 2513.8
$ws.Range("H74").Value = 2646.125
$ws.Range("M74").Value = -1577.8
$ws.Range("N74").Value = -4738.6667
$ws.Range("L74").Value = 2866.6667
$ws.Range("M77").Value = -7889
$ws.Range("N77").Value = -23693.3335
$ws.Range("K77").Value = 12569
$ws.Range("I77").Value = 2513.8
$ws.Range("J77").Value = 2866.6667
$ws.Range("H77").Value = 2646.125
$ws.Range("L77").Value = 14333.3335
$ws.Range("H138").Value = 3569.4
$ws.Range("L138").Value = 21984.429
$ws.Range("M138").Value = 503.6154999999999
$ws.Range("N138").Value = -32264.429
$ws.Range("K138").Value = 4636.3845
$ws.Range("I138").Value = 1545.4615
$ws.Range("J138").Value = 7328.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J27").Value = 43004
$ws.Range("H27").Value = 43004
$ws.Range("N27").Value = -43372
$ws.Range("L27").Value = 43004
$ws.Range("M38").Value = -1372.6666
$ws.Range("K38").Value = 1839.6666
$ws.Range("I38").Value = 1839.6666
$ws.Range("H38").Value = 1839.6666
$ws.Range("I132").Value = 9789.357
$ws.Range("J132").Value = 23785.2
$ws.Range("H132").Value = 13472.474
$ws.Range("L132").Value = 71355.60000000001
$ws.Range("N132").Value = -76415.60000000001
$ws.Range("M132").Value = -26838.071
$ws.Range("K132").Value = 29368.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M7").Value = 72.90000000000001
$ws.Range("L7").Value = 74.304344
$ws.Range("N7").Value = -300.304344
$ws.Range("K7").Value = 40.1
$ws.Range("J7").Value = 74.304344
$ws.Range("I7").Value = 40.1
$ws.Range("H7").Value = 63.939392
$ws.Range("N22").Value = -890
$ws.Range("J22").Value = 190
$ws.Range("H22").Value = 496.5
$ws.Range("L22").Value = 190
$ws.Range("L38").Value = 7000
$ws.Range("M38").Value = -1423.5
$ws.Range("N38").Value = -7754
$ws.Range("K38").Value = 1800.5
$ws.Range("I38").Value = 1800.5
$ws.Range("J38").Value = 7000
$ws.Range("H38").Value = 2840.4
$ws.Range("K46").Value = 1800.5
$ws.Range("N46").Value = -7422
$ws.Range("I46").Value = 1800.5
$ws.Range("J46").Value = 7000
$ws.Range("H46").Value = 2840.4
$ws.Range("M46").Value = -1589.5
$ws.Range("L46").Value = 7000
$ws.Range("J94").Value = 1233.3334
$ws.Range("H94").Value = 1140
$ws.Range("L94").Value = 1233.3334
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -2135.3334
$ws.Range("K94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("I132").Value = 19233876
$ws.Range("J132").Value = 2493.0588
$ws.Range("H132").Value = 14495709
$ws.Range("L132").Value = 7479.176399999999
$ws.Range("N132").Value = -12539.1764
$ws.Range("M132").Value = -57699098
$ws.Range("K132").Value = 57701628

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J20").Value = 3000
$ws.Range("N20").Value = -9454
$ws.Range("H20").Value = 2375
$ws.Range("L20").Value = 9000
$ws.Range("I131").Value = 315.23254
$ws.Range("H131").Value = 639.5599999999999
$ws.Range("M131").Value = 4094.30238
$ws.Range("L131").Value = 2652.6843
$ws.Range("K131").Value = 945.6976199999999
$ws.Range("N131").Value = -12732.6843
$ws.Range("J131").Value = 884.2281

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 21.538462
$ws.Range("H2").Value = 31.882353
$ws.Range("M2").Value = 91.461538
$ws.Range("K2").Value = 21.538462
$ws.Range("M13").Value = -46.71428
$ws.Range("L13").Value = 800
$ws.Range("K13").Value = 185.71428
$ws.Range("N13").Value = -1078
$ws.Range("J13").Value = 800
$ws.Range("I13").Value = 185.71428
$ws.Range("H13").Value = 322.22223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L16").Value = 1246
$ws.Range("K16").Value = 999.5
$ws.Range("N16").Value = -1586
$ws.Range("J16").Value = 1246
$ws.Range("I16").Value = 999.5
$ws.Range("H16").Value = 1147.4
$ws.Range("M16").Value = -829.5
$ws.Range("M22").Value = -2037.3333
$ws.Range("N22").Value = -6088
$ws.Range("K22").Value = 2332.3333
$ws.Range("I22").Value = 2332.3333
$ws.Range("J22").Value = 5498
$ws.Range("H22").Value = 3598.6
$ws.Range("L22").Value = 5498
$ws.Range("J27").Value = 5498
$ws.Range("I27").Value = 2332.3333
$ws.Range("H27").Value = 3598.6
$ws.Range("N27").Value = -5712
$ws.Range("M27").Value = -2225.3333
$ws.Range("L27").Value = 5498
$ws.Range("K27").Value = 2332.3333
$ws.Range("K46").Value = 2654.2
$ws.Range("N46").Value = -6542.6665
$ws.Range("I46").Value = 2654.2
$ws.Range("J46").Value = 6166.6665
$ws.Range("H46").Value = 3971.375
$ws.Range("M46").Value = -2466.2
$ws.Range("L46").Value = 6166.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N47").ClearContents()
$ws.Range("J47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("L126").Value = 3250.7142
$ws.Range("M126").Value = 175.75
$ws.Range("N126").Value = -8190.7142
$ws.Range("K126").Value = 2294.25
$ws.Range("I126").Value = 764.75
$ws.Range("J126").Value = 1083.5714
$ws.Range("H126").Value = 882.2105
